$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 85
$prevRow = 84

# Row 84's Date / Weekday / Week columns already hold the same literal text
# values ("2025-02-23" / "Sunday" / "08") that the new row needs, and those
# cells are stored as plain text (no numeric/date auto-detection). Copying
# the whole row down preserves that text typing (and the numeric cell
# formatting for the count columns) instead of re-triggering Excel's
# smart type inference, which would otherwise turn "2025-02-23" into a date
# serial and "08" into the number 8.
$srcRange = "A" + $prevRow + ":T" + $prevRow
$dstRange = "A" + $row + ":T" + $row
$ws.Range($srcRange).Copy()
$ws.Range($dstRange).PasteSpecial()

$ws.Cells.Item($row, 2).Value = "23:33:12"
$ws.Cells.Item($row, 5).Value = 130465
$ws.Cells.Item($row, 6).Value = 141862
$ws.Cells.Item($row, 7).Value = 172541
$ws.Cells.Item($row, 8).Value = 158686
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 146707
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 193724
$ws.Cells.Item($row, 14).Value = 115590
$ws.Cells.Item($row, 15).Value = 46468
$ws.Cells.Item($row, 16).Value = 29360
$ws.Cells.Item($row, 17).Value = 68688
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 49221
$ws.Cells.Item($row, 20).Value = -1
